$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Insert two new rows at row 5, shifting the existing rows 5-7 down to 7-9
$ws.Rows.Item(5).Resize(2).Insert()

# New row 5: GONZALEZ / IXMATLAHUA / MIGUEL ANGEL - 4BEM
$ws.Cells.Item(5, 1).Value = 19330051920102
$ws.Cells.Item(5, 2).Value = "GONZALEZ"
$ws.Cells.Item(5, 3).Value = "IXMATLAHUA"
$ws.Cells.Item(5, 4).Value = "MIGUEL ANGEL"
$ws.Cells.Item(5, 5).Value = "FÍSICA I"
$ws.Cells.Item(5, 6).Value = "4BEM"
$ws.Cells.Item(5, 7).Value = 2

# New row 6: SARMIENTO / HERNANDEZ / ROMARIO ALDAIR - 4BEM
$ws.Cells.Item(6, 1).Value = 19330051920118
$ws.Cells.Item(6, 2).Value = "SARMIENTO"
$ws.Cells.Item(6, 3).Value = "HERNANDEZ"
$ws.Cells.Item(6, 4).Value = "ROMARIO ALDAIR"
$ws.Cells.Item(6, 5).Value = "FÍSICA I"
$ws.Cells.Item(6, 6).Value = "4BEM"
$ws.Cells.Item(6, 7).Value = 2
